$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the row above (A2:C2) into the new row 8,
# so the new cells reuse the existing cell styles (date / integer / wrap-text).
$ws.Range("A2:C2").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)

# Fill in the new journal entry.
$ws.Range("A8").Value = 43893
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "Clôture du sprint 1 et 2"

# Move the active selection, as recorded in the saved view state.
$ws.Range("C12").Select()
